$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "577.36") must keep
# their original Text cell type (matching the source inlineStr data) rather
# than being auto-converted to a Number by COM value assignment, so mark
# them as Text cells first. (Range.NumberFormat on a multi-area union only
# touches the first area in this host, so it is applied cell-by-cell.)
$textValueCells = @(
    'D5',
    'D6',
    'D9',
    'D10',
    'D12',
    'D16',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D28',
    'D29',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D41',
    'D42',
    'D44',
    'D45',
    'D49',
    'D50',
    'D51'
)
foreach ($ref in $textValueCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.629.32'
$ws.Range('E2').Value = '  +6.03%  '
$ws.Range('D3').Value = '3.395.38'
$ws.Range('E3').Value = '  +6.41%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '577.36'
$ws.Range('E5').Value = '  +7.70%  '
$ws.Range('D6').Value = '154.44'
$ws.Range('E6').Value = '  +6.44%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '3.398.31'
$ws.Range('E8').Value = '  +6.32%  '
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('D10').Value = '7.46'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E11').Value = '  +7.05%  '
$ws.Range('D12').Value = '0.438'
$ws.Range('E12').Value = '  +1.71%  '
$ws.Range('D13').Value = '3.978.55'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('E15').Value = '  +7.37%  '
$ws.Range('D16').Value = '27.10'
$ws.Range('E16').Value = '  +5.01%  '
$ws.Range('D17').Value = '63.636.56'
$ws.Range('E17').Value = '  +6.02%  '
$ws.Range('D18').Value = '3.400.71'
$ws.Range('E18').Value = '  +6.16%  '
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('E20').Value = '  +5.14%  '
$ws.Range('D21').Value = '8.44'
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('D22').Value = '390.37'
$ws.Range('E22').Value = '  +5.36%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '0.536'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('D25').Value = '71.06'
$ws.Range('E25').Value = '  +2.30%  '
$ws.Range('D26').Value = '9.59'
$ws.Range('E26').Value = '  +11.45%  '
$ws.Range('E27').Value = '  +18.63%  '
$ws.Range('D28').Value = '0.179'
$ws.Range('E28').Value = '  +6.16%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +7.87%  '
$ws.Range('D31').Value = '6.46'
$ws.Range('E31').Value = '  +5.74%  '
$ws.Range('D32').Value = '23.16'
$ws.Range('E32').Value = '  +2.96%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.32'
$ws.Range('E33').Value = '  +10.78%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '5.59'
$ws.Range('E34').Value = '  +5.92%  '
$ws.Range('D35').Value = '6.73'
$ws.Range('E35').Value = '  +2.76%  '
$ws.Range('D36').Value = '1.49'
$ws.Range('E36').Value = '  +9.40%  '
$ws.Range('D37').Value = '158.49'
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('D38').Value = '27.66'
$ws.Range('E38').Value = '  +5.29%  '
$ws.Range('E39').Value = '  +12.50%  '
$ws.Range('D40').Value = '2.892.39'
$ws.Range('E40').Value = '  +2.54%  '
$ws.Range('D41').Value = '0.0748'
$ws.Range('E41').Value = '  +6.27%  '
$ws.Range('D42').Value = '0.0325'
$ws.Range('E42').Value = '  +4.89%  '
$ws.Range('E43').Value = '  +6.28%  '
$ws.Range('D44').Value = '41.13'
$ws.Range('E44').Value = '  +4.19%  '
$ws.Range('D45').Value = '4.30'
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('E46').Value = '  +7.72%  '
$ws.Range('D47').Value = '3.440.73'
$ws.Range('E47').Value = '  +6.35%  '
$ws.Range('E48').Value = '  +7.02%  '
$ws.Range('D49').Value = '301.15'
$ws.Range('E49').Value = '  +14.05%  '
$ws.Range('D50').Value = '0.103'
$ws.Range('E50').Value = '  -0.46%  '
$ws.Range('D51').Value = '6.31'
$ws.Range('E51').Value = '  +2.69%  '
